{"js": "// Hands On Demos - Day 1.\n//\n// 1) The bullet \"Create a new Java Project called \"FullStackDay1\".\" used to be\n//    split across two runs with a \"_GoBack\" bookmark sitting in between (left\n//    over from the last place Word was edited). Re-merge it into a single run\n//    with no bookmark in the middle.\n// 2) Word always keeps exactly one \"_GoBack\" bookmark, tracking the most\n//    recent edit location. Re-insert it (collapsed) on the empty paragraph\n//    just above the closing \"****\" line at the end of the document.\n\nconst body = context.document.body;\n\n// --- Step 1: collapse the \"FullStackDay1\" run split back into one run ----\nconst hits = body.search(\"FullStackDay1\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const hitRange = hits.items[0];\n  const para = hitRange.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  // Re-writing the whole paragraph's text collapses every run inside it\n  // (and any bookmark boundaries that fell inside the old run split) into a\n  // single run carrying the paragraph's run formatting.\n  para.insertText(para.text, \"Replace\");\n  await context.sync();\n}\n\n// --- Step 2: move the \"_GoBack\" bookmark to the end of the document -------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst targetParagraph = lastParagraph.getPrevious();\nconst targetRange = targetParagraph.getRange();\ntargetRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Hands On Demos - Day 1.\n#\n# 1) The bullet \"Create a new Java Project called \"FullStackDay1\".\" used to be\n#    split across two runs with a \"_GoBack\" bookmark sitting in between (left\n#    over from the last place Word was edited). Re-merge it into a single run\n#    with no bookmark in the middle.\n# 2) Word always keeps exactly one \"_GoBack\" bookmark, tracking the most\n#    recent edit location. Re-insert it (collapsed) at the empty paragraph\n#    just above the closing \"****\" line at the end of the document.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: collapse the \"FullStackDay1\" run split back into one run ------\n# A Find/Replace across the whole story re-writes the matched span as a\n# single run and swallows any bookmark boundary that used to sit inside it.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"FullStackDay1\", $false, $false, $false, $false, $false, $true, 1, $false, \"FullStackDay1\", 2) | Out-Null\n\n# --- Step 2: move the \"_GoBack\" bookmark to the end of the document --------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$lastParagraph = $d.Paragraphs.Last\n$targetParagraph = $lastParagraph.Previous()\n$d.Bookmarks.Add(\"_GoBack\", $targetParagraph.Range)\n"}
